$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55 - this shifts the existing rows 55..75 down to 56..76
# (matching the "most recent report pushed on top, older rows shift down" pattern
# seen across the rest of the sheet).
$ws.Rows.Item(55).Insert()

# Populate the new row 55 with this week's reading for the same market/product
# series (same Mercado/Region/Producto/Categoria/Variedad/Calidad/Unidad/Origen
# as the row that is now directly below it), just a new Fecha + price reading.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 44992
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100108
$ws.Range("H55").Value = "Tropicales y subtropicales"
$ws.Range("I55").Value = 100108003
$ws.Range("J55").Value = "Maracuyá"
$ws.Range("K55").Value = "Sin especificar"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 20
$ws.Range("N55").Value = 60000
$ws.Range("O55").Value = 60000
$ws.Range("P55").Value = 60000
$ws.Range("Q55").Value = "$/caja 18 kilos"
$ws.Range("R55").Value = "Región de Arica y Parinacota"
$ws.Range("S55").Value = 3333
$ws.Range("T55").Value = 18
